$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the annotation scores for columns E:J (Clear, Assertive, Cautious,
# Optimistic, Specific, Relevant) across rows 2-29.
$values = @(
    @(2,2,1,2,1,2),
    @(1,2,1,1,1,2),
    @(2,2,1,1,2,2),
    @(2,1,1,1,1,2),
    @(2,2,1,1,2,2),
    @(2,2,1,2,1,2),
    @(2,0,1,1,2,2),
    @(2,2,2,2,2,2),
    @(2,2,2,2,2,2),
    @(2,2,2,2,2,2),
    @(2,2,2,2,2,2),
    @(2,1,1,1,1,2),
    @(2,2,1,1,1,2),
    @(2,1,1,1,2,2),
    @(2,2,1,1,2,2),
    @(2,0,1,0,2,2),
    @(2,2,1,1,2,2),
    @(2,2,1,1,2,2),
    @(2,2,1,2,2,2),
    @(2,2,1,2,2,2),
    @(2,2,2,1,1,2),
    @(2,2,1,1,2,2),
    @(2,1,2,1,1,2),
    @(2,2,2,1,1,2),
    @(2,2,1,1,2,2),
    @(2,2,1,1,2,2),
    @(2,2,1,2,2,2),
    @(2,2,2,2,2,2)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, 5 + $j).Value2 = $rowVals[$j]
    }
}

# Update the sheet view: zoom to 85%, freeze the header row, scroll so row 26
# is at the top of the frozen pane, and select H29.
$sheetView = $ws.Application.ActiveWindow
$ws.Select()
$excel.ActiveWindow.Zoom = 85

$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$excel.ActiveWindow.ScrollRow = 26
$ws.Range("H29").Select()

Write-Output "done"
